# Refresh the crypto price/volume snapshot to match the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '65.901.02'
$ws.Range("E2").Value = '  -1.21%  '

# Row 3
$ws.Range("D3").Value = '3.423.66'
$ws.Range("E3").Value = '  -1.61%  '

# Row 5
$ws.Range("D5").Value = '''597.40'
$ws.Range("E5").Value = '  -0.91%  '

# Row 6
$ws.Range("D6").Value = '''142.89'
$ws.Range("E6").Value = '  -3.05%  '

# Row 7
$ws.Range("D7").Value = '3.424.31'
$ws.Range("E7").Value = '  -1.47%  '

# Row 8
$ws.Range("E8").Value = '  -0.15%  '

# Row 9
$ws.Range("E9").Value = '  -2.30%  '

# Row 10
$ws.Range("D10").Value = '''8.06'
$ws.Range("E10").Value = '  +6.87%  '

# Row 11
$ws.Range("E11").Value = '  -5.19%  '

# Row 12
$ws.Range("D12").Value = '''0.406'
$ws.Range("E12").Value = '  -3.88%  '

# Row 13
$ws.Range("D13").Value = '4.000.97'
$ws.Range("E13").Value = '  -1.68%  '

# Row 14
$ws.Range("D14").Value = '''0.0000201'
$ws.Range("E14").Value = '  -5.80%  '

# Row 15
$ws.Range("D15").Value = '''29.71'
$ws.Range("E15").Value = '  -5.39%  '

# Row 16: WrappedEther
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.418.59'
$ws.Range("E16").Value = '  -1.48%  '

# Row 17: TRON
$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value = '''0.116'
$ws.Range("E17").Value = '  -0.72%  '

# Row 18
$ws.Range("D18").Value = '65.822.31'
$ws.Range("E18").Value = '  -1.45%  '

# Row 19
$ws.Range("D19").Value = '''10.40'
$ws.Range("E19").Value = '  +3.53%  '

# Row 20
$ws.Range("D20").Value = '''6.14'
$ws.Range("E20").Value = '  -4.46%  '

# Row 21
$ws.Range("D21").Value = '''14.61'
$ws.Range("E21").Value = '  -4.73%  '

# Row 22
$ws.Range("D22").Value = '''416.51'
$ws.Range("E22").Value = '  -4.95%  '

# Row 23
$ws.Range("D23").Value = '''0.579'
$ws.Range("E23").Value = '  -4.71%  '

# Row 24
$ws.Range("D24").Value = '''77.37'
$ws.Range("E24").Value = '  -2.60%  '

# Row 25
$ws.Range("E25").Value = '  +0.12%  '

# Row 26: PEPE
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").Value = '''0.0000110'
$ws.Range("E26").Value = '  -8.17%  '

# Row 27: InternetComputer(DFINITY)
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '''9.28'
$ws.Range("E27").Value = '  -5.24%  '

# Row 28: RenderToken
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").Value = '''7.86'
$ws.Range("E28").Value = '  -6.21%  '

# Row 29: PancakeSwap
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''2.42'
$ws.Range("E29").Value = '  -2.40%  '

# Row 30: Binance-PegBSC-USD
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -0.15%  '

# Row 31: Kaspa
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '''0.160'
$ws.Range("E31").Value = '  -4.26%  '

# Row 32: Fetch.AI
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = '''1.46'
$ws.Range("E32").Value = '  -7.99%  '

# Row 33: EthereumClassic
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '''24.65'
$ws.Range("E33").Value = '  -2.74%  '

# Row 34: RenzoRestakedETH
$ws.Range("B34").Value = 'RenzoRestakedETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D34").Value = '3.418.25'
$ws.Range("E34").Value = '  -1.53%  '

# Row 35: USDe
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  -0.10%  '

# Row 36: ImmutableX
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '''1.69'
$ws.Range("E36").Value = '  -6.15%  '

# Row 37: NEARProtocol
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").Value = '''5.50'
$ws.Range("E37").Value = '  -8.63%  '

# Row 38: Aptos
$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D38").Value = '''7.55'
$ws.Range("E38").Value = '  -4.51%  '

# Row 39: FirstDigitalUSD
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '''0.998'
$ws.Range("E39").Value = '  -0.15%  '

# Row 40: Monero
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '''168.75'
$ws.Range("E40").Value = '  -4.46%  '

# Row 41: Hedera
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '''0.0855'
$ws.Range("E41").Value = '  -3.12%  '

# Row 42: Mantle
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").Value = '''0.875'
$ws.Range("E42").Value = '  -1.87%  '

# Row 43: Filecoin
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = '''5.05'
$ws.Range("E43").Value = '  -6.83%  '

# Row 44: Stacks
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '''1.90'
$ws.Range("E44").Value = '  -10.31%  '

# Row 45: OKB
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '''45.46'
$ws.Range("E45").Value = '  -1.93%  '

# Row 46: InjectiveProtocol
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value = '''26.33'
$ws.Range("E46").Value = '  -8.77%  '

# Row 47: ONDO
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = '''1.18'
$ws.Range("E47").Value = '  -3.82%  '

# Row 48: Cosmos
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '''7.07'
$ws.Range("E48").Value = '  -5.04%  '

# Row 49: dogwifhat
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '''2.27'
$ws.Range("E49").Value = '  -6.61%  '

# Row 50: SuiNetwork
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").Value = '''0.920'
$ws.Range("E50").Value = '  -5.98%  '

# Row 51: TheGraph
$ws.Range("B51").Value = 'TheGraph'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D51").Value = '''0.232'
$ws.Range("E51").Value = '  -5.45%  '
